$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "56.249.53"
$ws.Range("E2").Value = "  -1.44%  "

# Row 3
$ws.Range("D3").Value = "2.368.84"
$ws.Range("E3").Value = "  -1.21%  "

# Row 4
$ws.Range("E4").Value = "  +0.16%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "506.00"
$ws.Range("E5").Value = "  +0.17%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "129.64"
$ws.Range("E6").Value = "  -2.50%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  +0.14%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.544"
$ws.Range("E8").Value = "  -1.76%  "

# Row 9
$ws.Range("D9").Value = "2.376.19"
$ws.Range("E9").Value = "  -1.43%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0983"
$ws.Range("E10").Value = "  +1.33%  "

# Row 11
$ws.Range("E11").Value = "  -0.08%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.85"
$ws.Range("E12").Value = "  +5.94%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.328"
$ws.Range("E13").Value = "  +2.00%  "

# Row 14
$ws.Range("D14").Value = "2.791.97"
$ws.Range("E14").Value = "  -1.20%  "

# Row 15
$ws.Range("D15").Value = "56.360.74"
$ws.Range("E15").Value = "  -1.05%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.46"
$ws.Range("E16").Value = "  -2.02%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000132"
$ws.Range("E17").Value = "  -1.16%  "

# Row 18
$ws.Range("D18").Value = "2.374.24"
$ws.Range("E18").Value = "  -1.02%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.01"
$ws.Range("E19").Value = "  -2.10%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.04"
$ws.Range("E20").Value = "  -0.02%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "308.67"
$ws.Range("E21").Value = "  -0.37%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.28"
$ws.Range("E22").Value = "  -0.62%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  +0.17%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "66.04"
$ws.Range("E24").Value = "  +1.65%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.998"
$ws.Range("E25").Value = "  -0.09%  "

# Row 26
$ws.Range("E26").Value = "  -1.50%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.148"
$ws.Range("E27").Value = "  -2.85%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.20"
$ws.Range("E28").Value = "  -3.07%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "173.17"
$ws.Range("E29").Value = "  -0.03%  "

# Row 30
$ws.Range("D30").Value = "0.0₃0709"
$ws.Range("E30").Value = "  -1.97%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.64"
$ws.Range("E31").Value = "  -2.02%  "

# Row 32
$ws.Range("B32").Value = "USDe"
$ws.Range("C32").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.999"
$ws.Range("E32").Value = "  +0.02%  "

# Row 33
$ws.Range("B33").Value = "Aptos"
$ws.Range("C33").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.82"
$ws.Range("E33").Value = "  -2.07%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.997"
$ws.Range("E34").Value = "  +0.13%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.06"
$ws.Range("E35").Value = "  -5.35%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.59"
$ws.Range("E36").Value = "  -1.94%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.18"
$ws.Range("E37").Value = "  -1.48%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.68"
$ws.Range("E38").Value = "  -4.14%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.820"
$ws.Range("E39").Value = "  +1.34%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "36.31"
$ws.Range("E40").Value = "  -1.16%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.38"
$ws.Range("E41").Value = "  -4.19%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.37"
$ws.Range("E42").Value = "  +0.49%  "

# Row 43
$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "125.82"
$ws.Range("E43").Value = "  -5.25%  "

# Row 44
$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.79"
$ws.Range("E44").Value = "  -4.56%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.562"
$ws.Range("E45").Value = "  -0.87%  "

# Row 46
$ws.Range("E46").Value = "  -1.37%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "238.54"
$ws.Range("E47").Value = "  -5.63%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0481"
$ws.Range("E48").Value = "  -1.75%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0206"
$ws.Range("E49").Value = "  -1.76%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "16.91"
$ws.Range("E50").Value = "  -1.85%  "

# Row 51
$ws.Range("E51").Value = "  -0.11%  "
